$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 5 (existing rows 5-10 shift down to 7-12)
$ws.Rows.Item(5).Resize(2).Insert()

# --- New row 5 ---
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 45014
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100107
$ws.Cells.Item(5, 8).Value = "Otros"
$ws.Cells.Item(5, 9).Value = 100107011
$ws.Cells.Item(5, 10).Value = "Tuna"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 50
$ws.Cells.Item(5, 14).Value = 13000
$ws.Cells.Item(5, 15).Value = 14000
$ws.Cells.Item(5, 16).Value = 13600
$ws.Cells.Item(5, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(5, 19).Value = 756
$ws.Cells.Item(5, 20).Value = 18

# --- New row 6 ---
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value = "Bíobío"
$ws.Cells.Item(6, 4).Value = 45014
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100107
$ws.Cells.Item(6, 8).Value = "Otros"
$ws.Cells.Item(6, 9).Value = 100107011
$ws.Cells.Item(6, 10).Value = "Tuna"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 20
$ws.Cells.Item(6, 14).Value = 10000
$ws.Cells.Item(6, 15).Value = 10000
$ws.Cells.Item(6, 16).Value = 10000
$ws.Cells.Item(6, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(6, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(6, 19).Value = 556
$ws.Cells.Item(6, 20).Value = 18
